$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 588, shifting existing rows 588-678 down to 589-679
$ws.Rows.Item(588).Insert()

# Populate the newly inserted row 588 with the new data record
$ws.Range("A588").Value = 11
$ws.Range("B588").Value = "Vega Monumental Concepción"
$ws.Range("C588").Value = "Bíobío"
$ws.Range("D588").Value = 44951
$ws.Range("E588").Value = 8
$ws.Range("F588").Value = "Fruta"
$ws.Range("G588").Value = 100102
$ws.Range("H588").Value = "Cítricos"
$ws.Range("I588").Value = 100102003
$ws.Range("J588").Value = "Limón"
$ws.Range("K588").Value = "Sin especificar"
$ws.Range("L588").Value = "1a amarillo"
$ws.Range("M588").Value = 170
$ws.Range("N588").Value = 12000
$ws.Range("O588").Value = 13000
$ws.Range("P588").Value = 12529
$ws.Range("Q588").Value = "$/malla 16 kilos"
$ws.Range("R588").Value = "Región de O'Higgins"
$ws.Range("S588").Value = 783
$ws.Range("T588").Value = 16
